$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated data values for columns A and B, rows 1-32
$values = @(
    @(-0.11638189569927704, 0.11633018141579754),
    @(-0.10461316739699278, 0.10444010341530152),
    @(-0.054737981802581714, 0.054603968806803849),
    @(-0.046603968846486765, 0.046133691450094716),
    @(-0.043133691469508406, 0.041530703103680899),
    @(-0.020299156455047651, 0.020069088019043946),
    @(-0.010069088072014676, 0.010024913038439554),
    @(-0.031915763476213499, 0.031638670715068429),
    @(-0.029638670736798378, 0.029405882915328974),
    @(-0.027405882939550708, 0.027390212493555666),
    @(-0.024390212522400034, 0.024363904832706673),
    @(-0.020863904864350413, 0.020670400350225826),
    @(-0.017170400384466156, 0.01708218901302061),
    @(-0.0090821890673336014, 0.0090533146878879478),
    @(-0.0080533147125620985, 0.0080346621295568355),
    @(-0.006034662158983739, 0.0060035529045565639),
    @(-0.0040035529346242882, 0.0039999999611506354),
    @(-0.0057170140451709983, 0.0056924942721465754),
    @(-0.0016924942899101403, 0.001535701873168982),
    @(-0.0080168288632442852, 0.0080056877501295531),
    @(-0.0040056877684602199, 0.0039999999815121257),
    @(-0.045702121633951265, 0.04549171451490075),
    @(-0.040491714541782464, 0.040097604659772301),
    @(-0.020097604753567921, 0.019999999904984911),
    @(-0.019231546666873101, 0.019158966263109534),
    @(-0.016658966283753074, 0.016568424314167274),
    @(-0.01406842433512967, 0.01355085246450205),
    @(-0.08897840267614221, 0.088350222687030744),
    @(-0.081350222734310584, 0.081168032946500546),
    @(-0.021168033224450156, 0.021023188721352781),
    @(-0.01402318877334352, 0.014001313373274371),
    @(-0.0040013134384793148, 0.003999999960491607)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
}

# Column B width changed (stored width 15.42578125 -> 14.7109375 character units)
$ws.Columns.Item(2).ColumnWidth = 13.75
